# Regenerate the "K" column (column G) of the save-data sheet.
#
# The save-data sheet used to store a raw "Strike#" count in column G.
# This regenerates that column from the (re-simulated) s_vals, writing the
# new "K" statistic per saved game row. The new values below are the
# freshly computed results for each row (2-71) of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row (A column index, 1-based Excel row) -> new K value
$kValues = [ordered]@{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 3
    13 = 0
    14 = 3
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 3
    24 = 2
    25 = 1
    26 = 1
    27 = 3
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 3
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 3
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 0
    46 = 2
    47 = 0
    48 = 1
    49 = 0
    50 = 2
    51 = 1
    52 = 0
    53 = 2
    54 = 1
    55 = 1
    56 = 3
    57 = 0
    58 = 2
    59 = 0
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 2
    65 = 2
    66 = 2
    67 = 3
    70 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
